$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.495.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.37%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.415.68'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.70%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '610.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.47%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.83%  '

# Row 7
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.595'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.48%  '

# Row 8
$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.409.74'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.70%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.02%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.192'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.67%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.87'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.05%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.557'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.07%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '43.68'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.46%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000267'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.14%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.969.50'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.52%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.06'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.75%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.422.08'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.36%  '

# Row 18
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.690.45'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.14%  '

# Row 19
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '574.35'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.30%  '

# Row 20
$ws.Range("E20").Value = '  +0.48%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.96'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.64%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.838'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.06%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.93'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.59%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '94.54'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.56%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.74%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.60'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.77%  '

# Row 27
$ws.Range("E27").Value = '  +0.02%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.38'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.93%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.37'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.06%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.47'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.54%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.75'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.55%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.42%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.75'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.02%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.49'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.19%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '582.35'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.99%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.39'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.84%  '

# Row 37
$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.46%  '

# Row 38
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0946'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.82%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '55.95'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.11%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0461'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.21%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.140'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.50%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.02'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -15.36%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.207.22'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.24%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₃0664'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -10.40%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '31.02'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.40%  '

# Row 46
$ws.Range("E46").Value = '  -5.29%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.290'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.36%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.35'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.85%  '

# Row 49
$ws.Range("E49").Value = '  -3.59%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.05'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.72%  '

# Row 51
$ws.Range("E51").Value = '  -0.02%  '
